$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's text value while preserving its exact string
# representation (Excel otherwise auto-converts numeric-looking text to
# numbers / strips trailing zeros) and without leaving a stray explicit
# cell style behind.
function Set-TextCell($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextCell $ws.Range("D2") "26.130.41"
Set-TextCell $ws.Range("E2") "  +1.37%  "

Set-TextCell $ws.Range("D3") "1.641.38"
Set-TextCell $ws.Range("E3") "  +0.41%  "

Set-TextCell $ws.Range("D4") "1.00"
Set-TextCell $ws.Range("E4") "  -0.08%  "

Set-TextCell $ws.Range("D5") "216.47"
Set-TextCell $ws.Range("E5") "  +0.30%  "

Set-TextCell $ws.Range("E6") "  +0.73%  "

Set-TextCell $ws.Range("E7") "  -0.12%  "

Set-TextCell $ws.Range("D8") "0.257"
Set-TextCell $ws.Range("E8") "  +0.56%  "

Set-TextCell $ws.Range("E9") "  +0.27%  "

Set-TextCell $ws.Range("E10") "  +0.59%  "

Set-TextCell $ws.Range("D11") "0.0791"
Set-TextCell $ws.Range("E11") "  -0.14%  "

Set-TextCell $ws.Range("D12") "1.868.55"
Set-TextCell $ws.Range("E12") "  +0.44%  "

Set-TextCell $ws.Range("E13") "  +0.66%  "

Set-TextCell $ws.Range("D14") "1.624.07"
Set-TextCell $ws.Range("E14") "  -0.67%  "

Set-TextCell $ws.Range("E15") "  -2.79%  "

Set-TextCell $ws.Range("E16") "  -0.10%  "

Set-TextCell $ws.Range("D17") "63.24"
Set-TextCell $ws.Range("E17") "  +0.03%  "

Set-TextCell $ws.Range("D18") "26.141.29"
Set-TextCell $ws.Range("E18") "  +1.33%  "

Set-TextCell $ws.Range("D19") "1.00"
Set-TextCell $ws.Range("E19") "  -0.10%  "

Set-TextCell $ws.Range("D20") "194.99"
Set-TextCell $ws.Range("E20") "  +1.19%  "

Set-TextCell $ws.Range("E21") "  -0.84%  "

Set-TextCell $ws.Range("D22") "10.02"
Set-TextCell $ws.Range("E22") "  +0.43%  "

Set-TextCell $ws.Range("E23") "  -0.38%  "

Set-TextCell $ws.Range("E24") "  -2.21%  "

Set-TextCell $ws.Range("E25") "  -0.13%  "

Set-TextCell $ws.Range("D26") "142.54"
Set-TextCell $ws.Range("E26") "  +0.15%  "

Set-TextCell $ws.Range("E27") "  +0.97%  "

Set-TextCell $ws.Range("E29") "  +0.56%  "

Set-TextCell $ws.Range("E30") "  +0.46%  "

Set-TextCell $ws.Range("E31") "  +1.77%  "

Set-TextCell $ws.Range("E32") "  +0.49%  "

Set-TextCell $ws.Range("D33") "3.24"
Set-TextCell $ws.Range("E33") "  +0.23%  "

Set-TextCell $ws.Range("E34") "  +1.22%  "

Set-TextCell $ws.Range("E35") "  +1.40%  "

Set-TextCell $ws.Range("D36") "0.911"
Set-TextCell $ws.Range("E36") "  +0.64%  "

Set-TextCell $ws.Range("D37") "1.136.10"
Set-TextCell $ws.Range("E37") "  +0.41%  "

Set-TextCell $ws.Range("E38") "  +1.28%  "

Set-TextCell $ws.Range("E39") "  -0.41%  "

Set-TextCell $ws.Range("E40") "  +1.13%  "

Set-TextCell $ws.Range("D41") "0.999"
Set-TextCell $ws.Range("E41") "  -0.23%  "

Set-TextCell $ws.Range("D42") "100.25"

Set-TextCell $ws.Range("E43") "  -1.25%  "

Set-TextCell $ws.Range("E44") "  -0.38%  "

Set-TextCell $ws.Range("D45") "1.777.96"
Set-TextCell $ws.Range("E45") "  +0.48%  "

Set-TextCell $ws.Range("D46") "0.0₆0112"
Set-TextCell $ws.Range("E46") "  -0.33%  "

Set-TextCell $ws.Range("D47") "56.78"
Set-TextCell $ws.Range("E47") "  +2.42%  "

Set-TextCell $ws.Range("E48") "  +3.85%  "

Set-TextCell $ws.Range("E49") "  +2.25%  "

Set-TextCell $ws.Range("B50") "Mantle"
Set-TextCell $ws.Range("C50") "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextCell $ws.Range("D50") "0.417"
Set-TextCell $ws.Range("E50") "  +0.07%  "

Set-TextCell $ws.Range("B51") "EnergySwap"
Set-TextCell $ws.Range("C51") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell $ws.Range("D51") "7.65"
Set-TextCell $ws.Range("E51") "  +2.86%  "
